# Update automatico via Actualizar 03-04-2021 22-41-18
#
# This mirrors the scheduled "Disponibilidad" checker run: the whole
# log rolls down one slot (rows 2-15 get the newest run's timestamp,
# rows 16-29 take what used to be the rows-2-15 timestamp, rows 30-43
# take what used to be the rows-16-29 timestamp), the status for the
# "Tomcat" row (row 4) flips back to "Disponible", and the check for
# that row additionally stamps an "Ultimo" (last-changed) timestamp
# into column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newest run timestamp, applied to rows 2-15 (column D).
$newTimestamp = 44259.94453363222

# Rows 16-29 inherit the timestamp that used to belong to rows 2-15.
$midTimestamp = 44250.35508177083

# Rows 30-43 inherit the timestamp that used to belong to rows 16-29.
$oldTimestamp = 44250.33371167824

# "Ultimo" (last-changed) timestamp stamped for row 4 only.
$ultimoTimestamp = 44259.9444875261

for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 4).Value = $newTimestamp
}

for ($row = 16; $row -le 29; $row++) {
    $ws.Cells.Item($row, 4).Value = $midTimestamp
}

for ($row = 30; $row -le 43; $row++) {
    $ws.Cells.Item($row, 4).Value = $oldTimestamp
}

# Row 4 ("Tomcat") flips from "No Disponible" back to "Disponible" and
# gets a new "Ultimo" timestamp in column E.
$ws.Cells.Item(4, 3).Value = "Disponible"
$ws.Cells.Item(4, 5).Value = $ultimoTimestamp
$ws.Cells.Item(4, 5).NumberFormat = $ws.Cells.Item(4, 4).NumberFormat

# Row 18 ("Tomcat" in the previous run's block) flips from "Disponible"
# to "No Disponible".
$ws.Cells.Item(18, 3).Value = "No Disponible"
